$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "78×48=3744"
$t.Cell(1, 2).Range.Text = "41×45=1845"
$t.Cell(1, 3).Range.Text = "52×78=4056"
$t.Cell(1, 4).Range.Text = "17×17=289"
$t.Cell(1, 5).Range.Text = "47×93=4371"
$t.Cell(2, 1).Range.Text = "90×64=5760"
$t.Cell(2, 2).Range.Text = "34×89=3026"
$t.Cell(2, 3).Range.Text = "37×96=3552"
$t.Cell(2, 4).Range.Text = "60×69=4140"
$t.Cell(2, 5).Range.Text = "97×33=3201"
$t.Cell(3, 1).Range.Text = "18×74=1332"
$t.Cell(3, 2).Range.Text = "12×82=984"
$t.Cell(3, 3).Range.Text = "78×15=1170"
$t.Cell(3, 4).Range.Text = "84×29=2436"
$t.Cell(3, 5).Range.Text = "63×53=3339"
$t.Cell(4, 1).Range.Text = "57×34=1938"
$t.Cell(4, 2).Range.Text = "12×44=528"
$t.Cell(4, 3).Range.Text = "37×53=1961"
$t.Cell(4, 4).Range.Text = "66×74=4884"
$t.Cell(4, 5).Range.Text = "94×93=8742"
$t.Cell(5, 1).Range.Text = "85×76=6460"
$t.Cell(5, 2).Range.Text = "72×70=5040"
$t.Cell(5, 3).Range.Text = "20×52=1040"
$t.Cell(5, 4).Range.Text = "35×29=1015"
$t.Cell(5, 5).Range.Text = "67×59=3953"
$t.Cell(6, 1).Range.Text = "12×22=264"
$t.Cell(6, 2).Range.Text = "84×48=4032"
$t.Cell(6, 3).Range.Text = "67×40=2680"
$t.Cell(6, 4).Range.Text = "98×46=4508"
$t.Cell(6, 5).Range.Text = "72×88=6336"
$t.Cell(7, 1).Range.Text = "67×87=5829"
$t.Cell(7, 2).Range.Text = "82×99=8118"
$t.Cell(7, 3).Range.Text = "68×81=5508"
$t.Cell(7, 4).Range.Text = "93×11=1023"
$t.Cell(7, 5).Range.Text = "64×39=2496"
$t.Cell(8, 1).Range.Text = "70×78=5460"
$t.Cell(8, 2).Range.Text = "70×18=1260"
$t.Cell(8, 3).Range.Text = "21×15=315"
$t.Cell(8, 4).Range.Text = "96×39=3744"
$t.Cell(8, 5).Range.Text = "62×69=4278"
$t.Cell(9, 1).Range.Text = "50×60=3000"
$t.Cell(9, 2).Range.Text = "59×21=1239"
$t.Cell(9, 3).Range.Text = "28×83=2324"
$t.Cell(9, 4).Range.Text = "15×27=405"
$t.Cell(9, 5).Range.Text = "94×63=5922"
$t.Cell(10, 1).Range.Text = "55×58=3190"
$t.Cell(10, 2).Range.Text = "15×41=615"
$t.Cell(10, 3).Range.Text = "22×31=682"
$t.Cell(10, 4).Range.Text = "48×61=2928"
$t.Cell(10, 5).Range.Text = "56×58=3248"
$t.Cell(11, 1).Range.Text = "21×63=1323"
$t.Cell(11, 2).Range.Text = "21×95=1995"
$t.Cell(11, 3).Range.Text = "17×65=1105"
$t.Cell(11, 4).Range.Text = "59×58=3422"
$t.Cell(11, 5).Range.Text = "11×77=847"
$t.Cell(12, 1).Range.Text = "44×17=748"
$t.Cell(12, 2).Range.Text = "23×84=1932"
$t.Cell(12, 3).Range.Text = "40×48=1920"
$t.Cell(12, 4).Range.Text = "86×59=5074"
$t.Cell(12, 5).Range.Text = "76×12=912"
$t.Cell(13, 1).Range.Text = "35×72=2520"
$t.Cell(13, 2).Range.Text = "74×40=2960"
$t.Cell(13, 3).Range.Text = "42×35=1470"
$t.Cell(13, 4).Range.Text = "55×100=5500"
$t.Cell(13, 5).Range.Text = "74×28=2072"
$t.Cell(14, 1).Range.Text = "28×58=1624"
$t.Cell(14, 2).Range.Text = "48×28=1344"
$t.Cell(14, 3).Range.Text = "63×89=5607"
$t.Cell(14, 4).Range.Text = "86×11=946"
$t.Cell(14, 5).Range.Text = "11×15=165"
$t.Cell(15, 1).Range.Text = "10×29=290"
$t.Cell(15, 2).Range.Text = "61×81=4941"
$t.Cell(15, 3).Range.Text = "67×67=4489"
$t.Cell(15, 4).Range.Text = "63×25=1575"
$t.Cell(15, 5).Range.Text = "76×78=5928"
$t.Cell(16, 1).Range.Text = "49×45=2205"
$t.Cell(16, 2).Range.Text = "69×57=3933"
$t.Cell(16, 3).Range.Text = "90×49=4410"
$t.Cell(16, 4).Range.Text = "52×26=1352"
$t.Cell(16, 5).Range.Text = "41×13=533"
$t.Cell(17, 1).Range.Text = "29×90=2610"
$t.Cell(17, 2).Range.Text = "28×19=532"
$t.Cell(17, 3).Range.Text = "95×65=6175"
$t.Cell(17, 4).Range.Text = "51×10=510"
$t.Cell(17, 5).Range.Text = "18×65=1170"
$t.Cell(18, 1).Range.Text = "87×28=2436"
$t.Cell(18, 2).Range.Text = "96×66=6336"
$t.Cell(18, 3).Range.Text = "24×81=1944"
$t.Cell(18, 4).Range.Text = "67×78=5226"
$t.Cell(18, 5).Range.Text = "67×66=4422"
$t.Cell(19, 1).Range.Text = "66×92=6072"
$t.Cell(19, 2).Range.Text = "58×43=2494"
$t.Cell(19, 3).Range.Text = "56×35=1960"
$t.Cell(19, 4).Range.Text = "31×22=682"
$t.Cell(19, 5).Range.Text = "45×48=2160"
$t.Cell(20, 1).Range.Text = "32×94=3008"
$t.Cell(20, 2).Range.Text = "11×31=341"
$t.Cell(20, 3).Range.Text = "57×22=1254"
$t.Cell(20, 4).Range.Text = "28×15=420"
$t.Cell(20, 5).Range.Text = "39×10=390"
